# Weekly update: insert two new daily-price rows at the top of the
# Maracuyá data block (row 16), pushing the existing rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 16 (existing rows 16.. shift down to 18..)
$ws.Rows("16:17").Insert()

# --- Row 16: new "Especial" quality entry for the latest date ---
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C16").Value = "Arica y Parinacota"
$ws.Range("D16").Value = 44676
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100108
$ws.Range("H16").Value = "Tropicales y subtropicales"
$ws.Range("I16").Value = 100108003
$ws.Range("J16").Value = "Maracuyá"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Especial"
$ws.Range("M16").Value = 40
$ws.Range("N16").Value = 24000
$ws.Range("O16").Value = 25000
$ws.Range("P16").Value = 24500
$ws.Range("Q16").Value = "$/caja 20 kilos"
$ws.Range("R16").Value = "Región de Arica y Parinacota"
$ws.Range("S16").Value = 1225
$ws.Range("T16").Value = 20

# --- Row 17: new "Primera" quality entry for the latest date ---
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C17").Value = "Arica y Parinacota"
$ws.Range("D17").Value = 44676
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100108
$ws.Range("H17").Value = "Tropicales y subtropicales"
$ws.Range("I17").Value = 100108003
$ws.Range("J17").Value = "Maracuyá"
$ws.Range("K17").Value = "Sin especificar"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 70
$ws.Range("N17").Value = 21000
$ws.Range("O17").Value = 22000
$ws.Range("P17").Value = 21500
$ws.Range("Q17").Value = "$/caja 20 kilos"
$ws.Range("R17").Value = "Región de Arica y Parinacota"
$ws.Range("S17").Value = 1075
$ws.Range("T17").Value = 20
